$d = $word.ActiveDocument

# --- Paragraph 1 edits ---
$p1 = $d.Paragraphs.Item(1)

# 1) Add a paragraph border (top/left/bottom/right, each with 5pt space-from-text,
#    no explicit line) to the first paragraph.
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# 2) Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.LeftIndent = 11.25

# 3) Remove the trailing " " run that followed the ID text. That run is the
#    single character immediately after the ID text within paragraph 1.
$full = $p1.Range
$idLen = ("**ID__AFFARS_5333_topic_11__ID**").Length
$spaceRange = $d.Range($full.Start + $idLen, $full.Start + $idLen + 1)
$spaceRange.Delete()

# 4) Update the ID text itself.
$d.Content.Find.Execute("**ID__AFFARS_5333_topic_11__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5333_215__ID**", 2) | Out-Null
